$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows - beach name, latitude, longitude
$beaches = @(
    @("Panama City Beach", 30.176591999999999, -85.805488999999994),
    @("Destin Beach", 30.393681000000001, -86.495659000000003),
    @("Amelia Island", 30.626570000000001, -81.460853999999998),
    @("St. Augustine Beach", 29.85202, -81.267319000000001),
    @("Cocoa Beach", 28.320221, -80.608870999999994),
    @("Clearwater beach", 27.924440000000001, -82.841003000000001),
    @("Vero Beach", 27.63888, -80.389076000000003),
    @("Longboat Key", 27.412541999999998, -82.658989000000005),
    @("Fort De Soto Park, Tierra Verde", 27.62509, -82.712639999999993),
    @("Pass-a-Grille Beach, St. Petersburg", 27.695070000000001, -82.735900000000001),
    @("Siesta Key, Sarasota", 27.266190000000002, -82.545649999999995),
    @("Blowing Rocks Preserve, Jupiter", 26.905529999999999, -80.132080000000002),
    @("Jupiter Dog Beach, Jupiter", 26.934139999999999, -80.099739999999997),
    @("Lighthouse Beach Park, Sanibel Island", 26.317830000000001, -80.083519999999993),
    @("Captiva Island", 26.519850000000002, -82.189903000000001),
    @("Bonita Beach Dog Beach, Fort Myers", 26.43695, -81.92165),
    @("South Beach, Miami", 40.688599000000004, -86.762710999999996),
    @("Cape Florida State Park, Key Biscayne", 25.672470000000001, -80.155410000000003),
    @("Smathers Beach", 29.84056, -81.270579999999995),
    @("Rest Beach", 29.262450000000001, -81.156270000000006)
)

$row = 3
foreach ($beach in $beaches) {
    $ws.Cells.Item($row, 1).Value = $beach[0]
    $ws.Cells.Item($row, 2).Value = $beach[1]
    $ws.Cells.Item($row, 3).Value = $beach[2]
    $row += 2
}

# Header row (row 2): Beachd / Lat / Long - added after the data so the
# shared-string table picks up indices 21/22/23 (matching the target order)
$ws.Range("A2").Value = "Beachd"
$ws.Range("B2").Value = "Lat"
$ws.Range("C2").Value = "Long"

# Column A width (target XML width is 33.7109375; the COM layer here quantizes
# ColumnWidth to 1/6-character steps, so 32.83 is the closest input that lands
# on the nearest achievable stored width, 33.666666666666664)
$ws.Columns.Item(1).ColumnWidth = 32.83

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1

# Selection matches the final state
$ws.Range("C41").Select()
